$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 370
$ws.Range("D370").Value = 45021
$ws.Range("I370").Value = "Primera"
$ws.Range("J370").Value = 500
$ws.Range("K370").Value = 1200
$ws.Range("L370").Value = 1200
$ws.Range("M370").Value = 1200
$ws.Range("O370").Value = "Región Metropolitana"
$ws.Range("P370").Value = 1200

# Row 371
$ws.Range("D371").Value = 45021
$ws.Range("I371").Value = "Segunda"
$ws.Range("J371").Value = 500
$ws.Range("K371").Value = 900
$ws.Range("L371").Value = 900
$ws.Range("M371").Value = 900
$ws.Range("O371").Value = "Región Metropolitana"
$ws.Range("P371").Value = 900

# Row 372
$ws.Range("I372").Value = "Extra"
$ws.Range("K372").Value = 800
$ws.Range("L372").Value = 800
$ws.Range("M372").Value = 800
$ws.Range("P372").Value = 800

# Row 373
$ws.Range("D373").Value = 44264
$ws.Range("I373").Value = "Primera"
$ws.Range("J373").Value = 400
$ws.Range("K373").Value = 600
$ws.Range("L373").Value = 600
$ws.Range("M373").Value = 600
$ws.Range("P373").Value = 600

# Row 374
$ws.Range("D374").Value = 44264
$ws.Range("I374").Value = "Segunda"
$ws.Range("J374").Value = 400
$ws.Range("K374").Value = 500
$ws.Range("L374").Value = 500
$ws.Range("M374").Value = 500
$ws.Range("P374").Value = 500

# Row 375
$ws.Range("I375").Value = "Extra"
$ws.Range("J375").Value = 800
$ws.Range("K375").Value = 1100
$ws.Range("L375").Value = 1100
$ws.Range("M375").Value = 1100
$ws.Range("P375").Value = 1100

# Row 376
$ws.Range("D376").Value = 44952
$ws.Range("H376").Value = "Tuna"
$ws.Range("J376").Value = 1000
$ws.Range("K376").Value = 800
$ws.Range("M376").Value = 800
$ws.Range("P376").Value = 800

# Row 377
$ws.Range("D377").Value = 44952
$ws.Range("H377").Value = "Tuna"
$ws.Range("J377").Value = 500

# Row 378
$ws.Range("H378").Value = "Calameño"

# Row 379
$ws.Range("H379").Value = "Calameño"

# Row 380
$ws.Range("D380").Value = 44551
$ws.Range("H380").Value = "Tuna"
$ws.Range("I380").Value = "Primera"
$ws.Range("J380").Value = 800
$ws.Range("K380").Value = 700
$ws.Range("L380").Value = 800
$ws.Range("M380").Value = 750
$ws.Range("P380").Value = 750

# Row 381
$ws.Range("D381").Value = 44551
$ws.Range("H381").Value = "Tuna"
$ws.Range("I381").Value = "Segunda"
$ws.Range("J381").Value = 400
$ws.Range("K381").Value = 600
$ws.Range("L381").Value = 600
$ws.Range("M381").Value = 600
$ws.Range("P381").Value = 600

# Row 382
$ws.Range("I382").Value = "Extra"
$ws.Range("K382").Value = 1400
$ws.Range("L382").Value = 1400
$ws.Range("M382").Value = 1400
$ws.Range("P382").Value = 1400

# Row 383
$ws.Range("H383").Value = "Calameño"
$ws.Range("I383").Value = "Primera"
$ws.Range("K383").Value = 900
$ws.Range("L383").Value = 900
$ws.Range("M383").Value = 900
$ws.Range("P383").Value = 900

# Row 384
$ws.Range("H384").Value = "Calameño"
$ws.Range("I384").Value = "Segunda"
$ws.Range("J384").Value = 500
$ws.Range("K384").Value = 700
$ws.Range("L384").Value = 700
$ws.Range("M384").Value = 700
$ws.Range("P384").Value = 700

# Row 385
$ws.Range("I385").Value = "Extra"
$ws.Range("J385").Value = 1000
$ws.Range("K385").Value = 1200
$ws.Range("L385").Value = 1200
$ws.Range("M385").Value = 1200
$ws.Range("P385").Value = 1200

# Row 386
$ws.Range("D386").Value = 44918
$ws.Range("H386").Value = "Tuna"
$ws.Range("J386").Value = 1000
$ws.Range("K386").Value = 900
$ws.Range("L386").Value = 900
$ws.Range("M386").Value = 900
$ws.Range("P386").Value = 900

# Row 387
$ws.Range("D387").Value = 44918
$ws.Range("H387").Value = "Tuna"
$ws.Range("K387").Value = 700
$ws.Range("L387").Value = 700
$ws.Range("M387").Value = 700
$ws.Range("P387").Value = 700

# Row 388
$ws.Range("H388").Value = "Calameño"

# Row 389
$ws.Range("H389").Value = "Calameño"

# Row 390
$ws.Range("D390").Value = 44217
$ws.Range("K390").Value = 800
$ws.Range("L390").Value = 800
$ws.Range("M390").Value = 800
$ws.Range("O390").Value = "Región de O'Higgins"
$ws.Range("P390").Value = 800

# Row 391 (rewritten)
$ws.Range("A391").Value = 11
$ws.Range("B391").Value = "Vega Monumental Concepción"
$ws.Range("C391").Value = "Bíobío"
$ws.Range("D391").Value = 44217
$ws.Range("E391").Value = 8
$ws.Range("F391").Value = 100112027
$ws.Range("G391").Value = "Melón"
$ws.Range("H391").Value = "Tuna"
$ws.Range("I391").Value = "Segunda"
$ws.Range("J391").Value = 500
$ws.Range("K391").Value = 600
$ws.Range("L391").Value = 600
$ws.Range("M391").Value = 600
$ws.Range("N391").Value = "`$/unidad"
$ws.Range("O391").Value = "Región de O'Higgins"
$ws.Range("P391").Value = 600
$ws.Range("Q391").Value = 1
$ws.Range("R391").Value = "Hortaliza"
$ws.Range("D391").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 392 (new)
$ws.Range("A392").Value = 11
$ws.Range("B392").Value = "Vega Monumental Concepción"
$ws.Range("C392").Value = "Bíobío"
$ws.Range("D392").Value = 45007
$ws.Range("E392").Value = 8
$ws.Range("F392").Value = 100112027
$ws.Range("G392").Value = "Melón"
$ws.Range("H392").Value = "Tuna"
$ws.Range("I392").Value = "Primera"
$ws.Range("J392").Value = 500
$ws.Range("K392").Value = 1000
$ws.Range("L392").Value = 1000
$ws.Range("M392").Value = 1000
$ws.Range("N392").Value = "`$/unidad"
$ws.Range("O392").Value = "Región Metropolitana"
$ws.Range("P392").Value = 1000
$ws.Range("Q392").Value = 1
$ws.Range("R392").Value = "Hortaliza"
$ws.Range("D392").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 393 (new)
$ws.Range("A393").Value = 11
$ws.Range("B393").Value = "Vega Monumental Concepción"
$ws.Range("C393").Value = "Bíobío"
$ws.Range("D393").Value = 45007
$ws.Range("E393").Value = 8
$ws.Range("F393").Value = 100112027
$ws.Range("G393").Value = "Melón"
$ws.Range("H393").Value = "Tuna"
$ws.Range("I393").Value = "Segunda"
$ws.Range("J393").Value = 500
$ws.Range("K393").Value = 800
$ws.Range("L393").Value = 800
$ws.Range("M393").Value = 800
$ws.Range("N393").Value = "`$/unidad"
$ws.Range("O393").Value = "Región Metropolitana"
$ws.Range("P393").Value = 800
$ws.Range("Q393").Value = 1
$ws.Range("R393").Value = "Hortaliza"
$ws.Range("D393").NumberFormat = "YYYY-MM-DD HH:MM:SS"
